$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the mortgage cost value (N1 "BILLS" / row for Mortagage is O1)
$ws.Range("O1").Value = 1633

# Update the active selection to match the recorded cursor position
$ws.Range("L19").Select()
